$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Park Data")

for ($r = 284; $r -le 330; $r++) {
    $ws.Cells.Item($r, 7).Value = 46.85230749999999
    $ws.Cells.Item($r, 8).Value = -121.7603229
}
